$d = $word.ActiveDocument

# Locate the paragraph that contains "LOB1046: Engenharia do Meio Ambiente (Requisito fraco)"
# and remove the three paragraphs that directly follow it:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "LOB1046: Engenharia do Meio Ambiente \(Requisito fraco\)") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $start = $target.Next().Range.Start
    $end = $target.Next().Next().Next().Next().Range.Start
    $r = $d.Range($start, $end)
    $r.Delete()
}
